$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# =====================================================================
# Content edits ("sierra leone master data"): translate fra -> eng,
# translate the French descriptions, and store is_active as literal
# text "TRUE" rather than a boolean.
# =====================================================================
$ws.Range("A2").Value = "eng"
$ws.Range("C2").Value = "Activated"
$ws.Range("A3").Value = "eng"
$ws.Range("C3").Value = "Deactivated"

$ws.Range("E1:E3").NumberFormat = "@"
$ws.Range("E2").Value = "'TRUE"
$ws.Range("E3").Value = "'TRUE"

# =====================================================================
# Restyle the table: plain Calibri body, bold Cambria header with a
# thin border, everything left aligned (vertical top on the header),
# and drop the old thick-bottom-border / fixed row height formatting.
# =====================================================================
$ws.Cells.ClearFormats()

$all = $ws.Range("A1:E3")
$all.Font.Name = "Calibri"
$all.Font.Size = 11
$all.Font.Color = 0
$all.HorizontalAlignment = -4131

$ws.Range("E1:E3").NumberFormat = "@"

$hdr = $ws.Range("A1:E1")
$hdr.Font.Name = "Cambria"
$hdr.Font.Bold = $true
$hdr.VerticalAlignment = -4160
$hdr.Borders.LineStyle = 1

$ws.Rows("1:3").AutoFit()

# A couple of extra (empty, word-wrapped) rows left below the table.
$ws.Range("C6").WrapText = $true
$ws.Range("C7").WrapText = $true

# =====================================================================
# Column widths
# =====================================================================
$ws.Columns("C").ColumnWidth = 21
$ws.Columns("D").ColumnWidth = 14.5
$ws.Columns("E").ColumnWidth = 7.6

# =====================================================================
# View / window state
# =====================================================================
$excel.Windows.Item(1).Zoom = 100
$ws.Range("B7").Select()

# =====================================================================
# Page setup
# =====================================================================
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36.85
$ws.PageSetup.FooterMargin = 36.85
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

Write-Host "Edit complete"
